# Adds a new "Local da Ocorrência:" paragraph at the very end of the
# document body (right before the sectPr), matching the formatting of the
# other body paragraphs (Arial, not bold, not italic, 12pt, justified).

$d = $word.ActiveDocument

# Use an existing paragraph with the exact same run formatting we need
# (Arial / b=0 / i=0 / sz=24) as a formatting donor so the emitted OOXML
# reproduces the explicit <w:b w:val="0"/><w:i w:val="0"/> markers exactly.
$donorPara = $d.Paragraphs.Item(17)

# Collapse to the very end of the document and insert a new paragraph
# containing the target text in one shot (carriage return + text).
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter([char]13 + "Local da Ocorrência:")

# Grab the freshly created last paragraph and copy the donor's character
# formatting onto its text (excludes the trailing paragraph mark).
$newPara = $d.Paragraphs.Last
$newRun = $newPara.Range
$newTextRange = $d.Range($newRun.Start, $newRun.End - 1)
$newTextRange.FormattedText = $donorPara.Range.FormattedText

# FormattedText pulled in the donor's text too, so re-seat the range off
# the (possibly moved) last paragraph and restore the intended wording.
$newPara = $d.Paragraphs.Last
$newRun = $newPara.Range
$newTextRange = $d.Range($newRun.Start, $newRun.End - 1)
$newTextRange.Text = "Local da Ocorrência:"

# Match the donor paragraph's justification (both/justify).
$newPara.Alignment = 3

Write-Host "Inserted 'Local da Ocorrência:' paragraph at end of document."
